$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: set a cell to an exact text value (so numeric/percent-looking
# strings like "308.12" or "0.27%" are not reinterpreted by Excel as
# numbers / percentages). Only touches the specific cells being changed.
function Set-TextCell($range, $value) {
    $r = $ws.Range($range)
    if ($r.NumberFormat -ne "@") {
        $r.NumberFormat = "@"
    }
    $r.Value = $value
}

Set-TextCell "D2" "308.12"
Set-TextCell "E2" "0.27%"
Set-TextCell "D3" "39.88"
Set-TextCell "E3" "1.78%"
Set-TextCell "D4" "5.149"
Set-TextCell "E4" "1.20%"
Set-TextCell "D5" "0.08101"
Set-TextCell "E5" "-0.63%"
Set-TextCell "D6" "1.930"
Set-TextCell "E6" "-2.05%"
Set-TextCell "D7" "8.155"
Set-TextCell "E7" "3.26%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell "D8" "4.225"
Set-TextCell "E8" "1.18%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D9" "0.9282"
Set-TextCell "E9" "-0.04%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D10" "0.1428"
Set-TextCell "E10" "0.92%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D11" "0.1919"
Set-TextCell "E11" "-1.40%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D12" "0.09041"
Set-TextCell "E12" "-2.14%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D13" "0.03495"
Set-TextCell "E13" "-0.14%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D14" "0.09814"
Set-TextCell "E14" "-0.56%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D15" "0.001411"
Set-TextCell "E15" "-0.14%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell "D16" "0.005902"
Set-TextCell "E16" "1.35%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D17" "3.952"
Set-TextCell "E17" "0.16%"
Set-TextCell "D18" "3.362"
Set-TextCell "E18" "-2.67%"
Set-TextCell "D19" "0.3426"
Set-TextCell "E19" "-0.78%"
Set-TextCell "D20" "0.1348"
Set-TextCell "E20" "3.39%"
Set-TextCell "D21" "4.643"
Set-TextCell "E21" "-3.48%"
Set-TextCell "D22" "0.2425"
Set-TextCell "E22" "-7.40%"
Set-TextCell "D23" "0.04374"
Set-TextCell "E23" "-2.37%"
Set-TextCell "D24" "0.001216"
Set-TextCell "E24" "-2.15%"
Set-TextCell "D25" "0.004336"
Set-TextCell "E25" "3.90%"
Set-TextCell "D26" "0.0001302"
Set-TextCell "E26" "0.03%"
Set-TextCell "D27" "0.0004007"
Set-TextCell "E27" "-9.91%"
Set-TextCell "D39" "0.02029"
Set-TextCell "E39" "-3.78%"
Set-TextCell "D40" "0.05048"
Set-TextCell "E40" "-1.96%"
Set-TextCell "D41" "0.007386"
Set-TextCell "E41" "-1.20%"
Set-TextCell "D42" "0.009771"
Set-TextCell "E42" "-3.58%"
Set-TextCell "E43" "-0.23%"
Set-TextCell "D44" "0.002134"
Set-TextCell "E44" "0.03%"
Set-TextCell "D45" "0.008707"
Set-TextCell "E45" "-10.04%"
Set-TextCell "D46" "0.00006353"
Set-TextCell "E46" "0.46%"
Set-TextCell "D47" "0.00000000751"
Set-TextCell "E47" "0.03%"
Set-TextCell "D48" "0.002866"
Set-TextCell "E49" "-18.83%"
Set-TextCell "D50" "0.00002104"
Set-TextCell "E50" "0.03%"
Set-TextCell "D51" "0.0002003"
Set-TextCell "E51" "0.03%"
